# Insert a new "jitTmax" column before the existing "stimT" column (H).
# This shifts the old H:O columns to I:P and adds the new jitter-onset
# max-time parameter column with a value of 500 for every trial row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at position H (column 8); existing H:O shift right to I:P.
$ws.Columns.Item(8).Insert()

# New header label for the inserted column.
$ws.Cells.Item(1, 8).Value = "jitTmax"

# Fill the new column's values (500) for each data row (rows 2-15).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = 500
}

# Match the resulting active selection shown in the saved workbook.
$ws.Range("H2").Select()
